$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage so that
# numeric-looking strings (e.g. "593.84") are not coerced into numbers
# and lose their original textual formatting/trailing zeros.
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Rows 2-44, 47-51: refreshed Price (D) and Volume(1h) (E) figures
Set-TextValue "D2" '62.964.23'
$ws.Range("E2").Value = '  +2.73%  '
Set-TextValue "D3" '3.036.66'
$ws.Range("E3").Value = '  +1.81%  '
$ws.Range("E4").Value = '  -0.07%  '
Set-TextValue "D5" '593.84'
$ws.Range("E5").Value = '  -1.02%  '
Set-TextValue "D6" '154.17'
$ws.Range("E6").Value = '  +7.46%  '
$ws.Range("E7").Value = '  -0.01%  '
Set-TextValue "D8" '3.033.02'
$ws.Range("E8").Value = '  +1.80%  '
Set-TextValue "D9" '0.517'
$ws.Range("E9").Value = '  -0.33%  '
Set-TextValue "D10" '6.86'
$ws.Range("E10").Value = '  +13.32%  '
$ws.Range("E11").Value = '  +3.99%  '
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("E13").Value = '  +3.11%  '
Set-TextValue "D14" '35.82'
$ws.Range("E14").Value = '  +4.60%  '
$ws.Range("E15").Value = '  +0.35%  '
Set-TextValue "D16" '3.537.32'
$ws.Range("E16").Value = '  +1.68%  '
Set-TextValue "D17" '7.09'
$ws.Range("E17").Value = '  +2.48%  '
Set-TextValue "D18" '62.896.37'
$ws.Range("E18").Value = '  +2.58%  '
Set-TextValue "D19" '3.034.31'
$ws.Range("E19").Value = '  +1.68%  '
Set-TextValue "D20" '454.01'
$ws.Range("E20").Value = '  +1.11%  '
Set-TextValue "D21" '14.27'
$ws.Range("E21").Value = '  +0.90%  '
Set-TextValue "D22" '0.698'
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("E23").Value = '  +2.85%  '
Set-TextValue "D24" '83.13'
$ws.Range("E24").Value = '  +1.50%  '
Set-TextValue "D25" '11.17'
$ws.Range("E25").Value = '  +6.58%  '
$ws.Range("E26").Value = '  +4.69%  '
Set-TextValue "D27" '12.45'
$ws.Range("E27").Value = '  +4.34%  '
Set-TextValue "D29" '7.48'
$ws.Range("E29").Value = '  +4.83%  '
Set-TextValue "D30" '2.26'
$ws.Range("E30").Value = '  +10.11%  '
$ws.Range("E32").Value = '  -0.09%  '
Set-TextValue "D33" '27.62'
$ws.Range("E33").Value = '  +1.72%  '
$ws.Range("E34").Value = '  +1.66%  '
Set-TextValue "D35" '0.0₃0860'
$ws.Range("E35").Value = '  +4.66%  '
$ws.Range("E36").Value = '  +2.25%  '
$ws.Range("E37").Value = '  +3.11%  '
Set-TextValue "D38" '3.21'
$ws.Range("E38").Value = '  +11.92%  '
$ws.Range("E39").Value = '  +8.05%  '
$ws.Range("E40").Value = '  +2.62%  '
Set-TextValue "D41" '50.43'
$ws.Range("E41").Value = '  +0.19%  '
Set-TextValue "D42" '9.14'
$ws.Range("E42").Value = '  +0.61%  '
Set-TextValue "D43" '0.307'
$ws.Range("E43").Value = '  +14.19%  '
Set-TextValue "D44" '43.97'
$ws.Range("E44").Value = '  +10.71%  '
Set-TextValue "D47" '2.723.40'
$ws.Range("E47").Value = '  +1.31%  '
Set-TextValue "D48" '133.46'
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +7.20%  '
Set-TextValue "D51" '24.86'
$ws.Range("E51").Value = '  +6.42%  '

# Rows 45 and 46: VeChain and Bittensor swapped ranking positions
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D45" '392.00'
$ws.Range("E45").Value = '  -1.32%  '

$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D46" '0.0362'
$ws.Range("E46").Value = '  +3.25%  '

